# Applies the data rotation described by the diff:
#   - Row 4 takes on the species/location data that used to be in Row 6,
#     plus Row 6's old A/Q/R (id, easting, northing) values.
#   - Row 5 keeps its own species/location data but takes on Row 4's old
#     A/Q/R (id, easting, northing) values.
#   - Row 6 takes on the species/location data that used to be in Row 4,
#     plus Row 5's old A/Q/R (id, easting, northing) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 4 : becomes the old "Barrfagerspindling" record (old row 6) ----
$ws.Range("A4").Value = 97022737
$ws.Range("B4").Value = 85254
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 249228
$ws.Range("F4").Value = "Barrfagerspindling"
$ws.Range("G4").Value = "Cortinarius piceae"
$ws.Range("H4").Value = "Frøslev, T.S.Jeppesen & Brandrud"
$ws.Range("P4").Value = "Abborrselmon, väst om kraftledning i NV, Jmt"
$ws.Range("Q4").Value = 542677.5135761717
$ws.Range("R4").Value = 6969417.407723305
$ws.Range("AI4").Value = "gammal kalkpåverkad barrskog i nederkant av sandig nipa"

# ---- Row 5 : same species data, new id / coordinates ----
$ws.Range("A5").Value = 97022725
$ws.Range("Q5").Value = 542935.3525043444
$ws.Range("R5").Value = 6969356.376306581

# ---- Row 6 : becomes the old "Lammticka" record (old row 4) ----
$ws.Range("A6").Value = 97022729
$ws.Range("B6").Value = 90130
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 1958
$ws.Range("F6").Value = "Lammticka"
$ws.Range("G6").Value = "Albatrellus subrubescens"
$ws.Range("H6").Value = "(Murrill) Pouzar"
$ws.Range("P6").Value = "Abborrselmon, öst om kraftledning i NV, Jmt"
$ws.Range("Q6").Value = 542861.7670970106
$ws.Range("R6").Value = 6969387.384421867
$ws.Range("AI6").Value = "gammal kalkpåverkad sandtallskog med graninslag i nederkant av nipa"

Write-Output "Applied row 4/5/6 rotation"
